$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 31.645482
$ws.Range("H2").Value = 63.290964
$ws.Range("I2").Value = 0.03555980726701226
$ws.Range("J2").Value = 0.0244846141215985
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 13.939895
$ws.Range("N2").Value = 27.87979
$ws.Range("O2").Value = 0.4802889349037177
$ws.Range("P2").Value = 0.4091514124048364
$ws.Range("Q2").Value = 441.13469630439
$ws.Range("R2").Value = 1764.53878521756
$ws.Range("S2").Value = 0.0170789819576548
$ws.Range("T2").Value = 0.01001791445003943

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 31.645482
$ws.Range("H3").Value = 63.290964
$ws.Range("I3").Value = 0.03555980726701226
$ws.Range("J3").Value = 0.0244846141215985
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5122363333333334
$ws.Range("N3").Value = 1.536709
$ws.Range("O3").Value = 0.01764872999083942
$ws.Range("P3").Value = 0.02255205859890708
$ws.Range("Q3").Value = 16.209965666246
$ws.Range("R3").Value = 97.25979399747601
$ws.Range("S3").Value = 0.0006275854369817887
$ws.Range("T3").Value = 0.0005521784524419172

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 31.645482
$ws.Range("H4").Value = 63.290964
$ws.Range("I4").Value = 0.03555980726701226
$ws.Range("J4").Value = 0.0244846141215985
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1702606666666667
$ws.Range("N4").Value = 0.5107820000000001
$ws.Range("O4").Value = 0.005866207331499287
$ws.Range("P4").Value = 0.007496009716391951
$ws.Range("Q4").Value = 5.387980862308002
$ws.Range("R4").Value = 32.327885173848
$ws.Range("S4").Value = 0.0002086012020964489
$ws.Range("T4").Value = 0.0001835369053576099

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 31.645482
$ws.Range("H5").Value = 63.290964
$ws.Range("I5").Value = 0.03555980726701226
$ws.Range("J5").Value = 0.0244846141215985
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.194159
$ws.Range("N5").Value = 0.5824769999999999
$ws.Range("O5").Value = 0.006689607010093757
$ws.Range("P5").Value = 0.008548173685789306
$ws.Range("Q5").Value = 6.144255139637999
$ws.Range("R5").Value = 36.865530837828
$ws.Range("S5").Value = 0.0002378811359709881
$ws.Range("T5").Value = 0.0002092987341409536

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 31.645482
$ws.Range("H6").Value = 63.290964
$ws.Range("I6").Value = 0.03555980726701226
$ws.Range("J6").Value = 0.0244846141215985
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.215910000000001
$ws.Range("N6").Value = 27.64773
$ws.Range("O6").Value = 0.317527470477254
$ws.Range("P6").Value = 0.4057458029378115
$ws.Range("Q6").Value = 291.64191401862
$ws.Range("R6").Value = 1749.85148411172
$ws.Range("S6").Value = 0.01129121565215308
$ws.Range("T6").Value = 0.009934529416390463

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 31.645482
$ws.Range("H7").Value = 63.290964
$ws.Range("I7").Value = 0.03555980726701226
$ws.Range("J7").Value = 0.0244846141215985
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 4.991516
$ws.Range("N7").Value = 9.983032
$ws.Range("O7").Value = 0.1719790502865958
$ws.Range("P7").Value = 0.1465065426562639
$ws.Range("Q7").Value = 157.958929730712
$ws.Range("R7").Value = 631.8357189228481
$ws.Range("S7").Value = 0.006115541882155156
$ws.Range("T7").Value = 0.003587156163228131

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 83.22744366666667
$ws.Range("H8").Value = 249.682331
$ws.Range("I8").Value = 0.09352209759714789
$ws.Range("J8").Value = 0.09659160077758068
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 13.939895
$ws.Range("N8").Value = 27.87979
$ws.Range("O8").Value = 0.4802889349037177
$ws.Range("P8").Value = 0.4091514124048364
$ws.Range("Q8").Value = 1160.181825831748
$ws.Range("R8").Value = 6961.09095499049
$ws.Range("S8").Value = 0.04491762864489569
$ws.Range("T8").Value = 0.03952058988459123

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 83.22744366666667
$ws.Range("H9").Value = 249.682331
$ws.Range("I9").Value = 0.09352209759714789
$ws.Range("J9").Value = 0.09659160077758068
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5122363333333334
$ws.Range("N9").Value = 1.536709
$ws.Range("O9").Value = 0.01764872999083942
$ws.Range("P9").Value = 0.02255205859890708
$ws.Range("Q9").Value = 42.6321205765199
$ws.Range("R9").Value = 383.689085188679
$ws.Range("S9").Value = 0.001650546248668995
$ws.Range("T9").Value = 0.002178339440898239

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 83.22744366666667
$ws.Range("H10").Value = 249.682331
$ws.Range("I10").Value = 0.09352209759714789
$ws.Range("J10").Value = 0.09659160077758068
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.1702606666666667
$ws.Range("N10").Value = 0.5107820000000001
$ws.Range("O10").Value = 0.005866207331499287
$ws.Range("P10").Value = 0.007496009716391951
$ws.Range("Q10").Value = 14.17036004364912
$ws.Range("R10").Value = 127.533240392842
$ws.Range("S10").Value = 0.0005486200145815808
$ws.Range("T10").Value = 0.0007240515779505972

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 83.22744366666667
$ws.Range("H11").Value = 249.682331
$ws.Range("I11").Value = 0.09352209759714789
$ws.Range("J11").Value = 0.09659160077758068
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.194159
$ws.Range("N11").Value = 0.5824769999999999
$ws.Range("O11").Value = 0.006689607010093757
$ws.Range("P11").Value = 0.008548173685789306
$ws.Range("Q11").Value = 16.15935723487633
$ws.Range("R11").Value = 145.434215113887
$ws.Range("S11").Value = 0.000625626079684553
$ws.Range("T11").Value = 0.0008256817800351811

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 83.22744366666667
$ws.Range("H12").Value = 249.682331
$ws.Range("I12").Value = 0.09352209759714789
$ws.Range("J12").Value = 0.09659160077758068
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 9.215910000000001
$ws.Range("N12").Value = 27.64773
$ws.Range("O12").Value = 0.317527470477254
$ws.Range("P12").Value = 0.4057458029378115
$ws.Range("Q12").Value = 767.0166303620701
$ws.Range("R12").Value = 6903.149673258631
$ws.Range("S12").Value = 0.02969583508374925
$ws.Range("T12").Value = 0.03919163661454801

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 83.22744366666667
$ws.Range("H13").Value = 249.682331
$ws.Range("I13").Value = 0.09352209759714789
$ws.Range("J13").Value = 0.09659160077758068
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 4.991516
$ws.Range("N13").Value = 9.983032
$ws.Range("O13").Value = 0.1719790502865958
$ws.Range("P13").Value = 0.1465065426562639
$ws.Range("Q13").Value = 415.4311167012654
$ws.Range("R13").Value = 2492.586700207592
$ws.Range("S13").Value = 0.01608384152556782
$ws.Range("T13").Value = 0.01415130147955743

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 524.5768889999999
$ws.Range("H14").Value = 1573.730667
$ws.Range("I14").Value = 0.5894633891046084
$ws.Range("J14").Value = 0.6088102578564109
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 13.939895
$ws.Range("N14").Value = 27.87979
$ws.Range("O14").Value = 0.4802889349037177
$ws.Range("P14").Value = 0.4091514124048364
$ws.Range("Q14").Value = 7312.546752086654
$ws.Range("R14").Value = 43875.28051251992
$ws.Range("S14").Value = 0.2831127433177881
$ws.Range("T14").Value = 0.2490955768885031

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 524.5768889999999
$ws.Range("H15").Value = 1573.730667
$ws.Range("I15").Value = 0.5894633891046084
$ws.Range("J15").Value = 0.6088102578564109
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.5122363333333334
$ws.Range("N15").Value = 1.536709
$ws.Range("O15").Value = 0.01764872999083942
$ws.Range("P15").Value = 0.02255205859890708
$ws.Range("Q15").Value = 268.707342172767
$ws.Range("R15").Value = 2418.366079554903
$ws.Range("S15").Value = 0.01040328019379235
$ws.Range("T15").Value = 0.01372992461079351

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 524.5768889999999
$ws.Range("H16").Value = 1573.730667
$ws.Range("I16").Value = 0.5894633891046084
$ws.Range("J16").Value = 0.6088102578564109
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1702606666666667
$ws.Range("N16").Value = 0.5107820000000001
$ws.Range("O16").Value = 0.005866207331499287
$ws.Range("P16").Value = 0.007496009716391951
$ws.Range("Q16").Value = 89.31481083906601
$ws.Range("R16").Value = 803.833297551594
$ws.Range("S16").Value = 0.003457914454815871
$ws.Range("T16").Value = 0.004563647608330745

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 524.5768889999999
$ws.Range("H17").Value = 1573.730667
$ws.Range("I17").Value = 0.5894633891046084
$ws.Range("J17").Value = 0.6088102578564109
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.194159
$ws.Range("N17").Value = 0.5824769999999999
$ws.Range("O17").Value = 0.006689607010093757
$ws.Range("P17").Value = 0.008548173685789306
$ws.Range("Q17").Value = 101.851324191351
$ws.Range("R17").Value = 916.6619177221587
$ws.Range("S17").Value = 0.003943278419947812
$ws.Range("T17").Value = 0.005204215825846774

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 524.5768889999999
$ws.Range("H18").Value = 1573.730667
$ws.Range("I18").Value = 0.5894633891046084
$ws.Range("J18").Value = 0.6088102578564109
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 9.215910000000001
$ws.Range("N18").Value = 27.64773
$ws.Range("O18").Value = 0.317527470477254
$ws.Range("P18").Value = 0.4057458029378115
$ws.Range("Q18").Value = 4834.45339710399
$ws.Range("R18").Value = 43510.08057393591
$ws.Range("S18").Value = 0.1871708188813357
$ws.Range("T18").Value = 0.2470222069107255

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 524.5768889999999
$ws.Range("H19").Value = 1573.730667
$ws.Range("I19").Value = 0.5894633891046084
$ws.Range("J19").Value = 0.6088102578564109
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 4.991516
$ws.Range("N19").Value = 9.983032
$ws.Range("O19").Value = 0.1719790502865958
$ws.Range("P19").Value = 0.1465065426562639
$ws.Range("Q19").Value = 2618.433934673724
$ws.Range("R19").Value = 15710.60360804234
$ws.Range("S19").Value = 0.1013753538369286
$ws.Range("T19").Value = 0.08919468601221126

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 181.4813383333333
$ws.Range("H20").Value = 544.4440149999999
$ws.Range("I20").Value = 0.2039293133121744
$ws.Range("J20").Value = 0.210622508737405
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 13.939895
$ws.Range("N20").Value = 27.87979
$ws.Range("O20").Value = 0.4802889349037177
$ws.Range("P20").Value = 0.4091514124048364
$ws.Range("Q20").Value = 2529.830800826141
$ws.Range("R20").Value = 15178.98480495685
$ws.Range("S20").Value = 0.09794499268635078
$ws.Range("T20").Value = 0.08617649693415923

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 181.4813383333333
$ws.Range("H21").Value = 544.4440149999999
$ws.Range("I21").Value = 0.2039293133121744
$ws.Range("J21").Value = 0.210622508737405
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.5122363333333334
$ws.Range("N21").Value = 1.536709
$ws.Range("O21").Value = 0.01764872999083942
$ws.Range("P21").Value = 0.02255205859890708
$ws.Range("Q21").Value = 92.96133531629278
$ws.Range("R21").Value = 836.652017846635
$ws.Range("S21").Value = 0.00359909338786386
$ws.Range("T21").Value = 0.004749971159294776

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 181.4813383333333
$ws.Range("H22").Value = 544.4440149999999
$ws.Range("I22").Value = 0.2039293133121744
$ws.Range("J22").Value = 0.210622508737405
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 0.1702606666666667
$ws.Range("N22").Value = 0.5107820000000001
$ws.Range("O22").Value = 0.005866207331499287
$ws.Range("P22").Value = 0.007496009716391951
$ws.Range("Q22").Value = 30.89913365219222
$ws.Range("R22").Value = 278.09220286973
$ws.Range("S22").Value = 0.001196291632859493
$ws.Range("T22").Value = 0.001578828371986436

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 181.4813383333333
$ws.Range("H23").Value = 544.4440149999999
$ws.Range("I23").Value = 0.2039293133121744
$ws.Range("J23").Value = 0.210622508737405
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.194159
$ws.Range("N23").Value = 0.5824769999999999
$ws.Range("O23").Value = 0.006689607010093757
$ws.Range("P23").Value = 0.008548173685789306
$ws.Range("Q23").Value = 35.23623516946166
$ws.Range("R23").Value = 317.1261165251549
$ws.Range("S23").Value = 0.001364206963896728
$ws.Range("T23").Value = 0.001800437786824013

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 181.4813383333333
$ws.Range("H24").Value = 544.4440149999999
$ws.Range("I24").Value = 0.2039293133121744
$ws.Range("J24").Value = 0.210622508737405
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 9.215910000000001
$ws.Range("N24").Value = 27.64773
$ws.Range("O24").Value = 0.317527470477254
$ws.Range("P24").Value = 0.4057458029378115
$ws.Range("Q24").Value = 1672.51568075955
$ws.Range("R24").Value = 15052.64112683595
$ws.Range("S24").Value = 0.06475315901217815
$ws.Range("T24").Value = 0.0854591989244346

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 181.4813383333333
$ws.Range("H25").Value = 544.4440149999999
$ws.Range("I25").Value = 0.2039293133121744
$ws.Range("J25").Value = 0.210622508737405
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 4.991516
$ws.Range("N25").Value = 9.983032
$ws.Range("O25").Value = 0.1719790502865958
$ws.Range("P25").Value = 0.1465065426562639
$ws.Range("Q25").Value = 905.8670039922465
$ws.Range("R25").Value = 5435.202023953479
$ws.Range("S25").Value = 0.03507156962902539
$ws.Range("T25").Value = 0.03085757556070593

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 15.79677433333333
$ws.Range("H26").Value = 47.390323
$ws.Range("I26").Value = 0.01775072507139627
$ws.Range("J26").Value = 0.01833332435500452
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 13.939895
$ws.Range("N26").Value = 27.87979
$ws.Range("O26").Value = 0.4802889349037177
$ws.Range("P26").Value = 0.4091514124048364
$ws.Range("Q26").Value = 220.2053755453617
$ws.Range("R26").Value = 1321.23225327217
$ws.Range("S26").Value = 0.008525476838309632
$ws.Range("T26").Value = 0.007501105553926084

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 15.79677433333333
$ws.Range("H27").Value = 47.390323
$ws.Range("I27").Value = 0.01775072507139627
$ws.Range("J27").Value = 0.01833332435500452
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 0.5122363333333334
$ws.Range("N27").Value = 1.536709
$ws.Range("O27").Value = 0.01764872999083942
$ws.Range("P27").Value = 0.02255205859890708
$ws.Range("Q27").Value = 8.091681763000778
$ws.Range("R27").Value = 72.825135867007
$ws.Range("S27").Value = 0.0003132777539266965
$ws.Range("T27").Value = 0.0004134542051668323

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 15.79677433333333
$ws.Range("H28").Value = 47.390323
$ws.Range("I28").Value = 0.01775072507139627
$ws.Range("J28").Value = 0.01833332435500452
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 0.1702606666666667
$ws.Range("N28").Value = 0.5107820000000001
$ws.Range("O28").Value = 0.005866207331499287
$ws.Range("P28").Value = 0.007496009716391951
$ws.Range("Q28").Value = 2.689569329176223
$ws.Range("R28").Value = 24.20612396258601
$ws.Range("S28").Value = 0.000104129433553253
$ws.Range("T28").Value = 0.0001374267774988791

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 15.79677433333333
$ws.Range("H29").Value = 47.390323
$ws.Range("I29").Value = 0.01775072507139627
$ws.Range("J29").Value = 0.01833332435500452
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 0.194159
$ws.Range("N29").Value = 0.5824769999999999
$ws.Range("O29").Value = 0.006689607010093757
$ws.Range("P29").Value = 0.008548173685789306
$ws.Range("Q29").Value = 3.067085907785666
$ws.Range("R29").Value = 27.603773170071
$ws.Range("S29").Value = 0.0001187453748718595
$ws.Range("T29").Value = 0.0001567164408244898

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 15.79677433333333
$ws.Range("H30").Value = 47.390323
$ws.Range("I30").Value = 0.01775072507139627
$ws.Range("J30").Value = 0.01833332435500452
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 9.215910000000001
$ws.Range("N30").Value = 27.64773
$ws.Range("O30").Value = 0.317527470477254
$ws.Range("P30").Value = 0.4057458029378115
$ws.Range("Q30").Value = 145.58165054631
$ws.Range("R30").Value = 1310.23485491679
$ws.Range("S30").Value = 0.005636342831057632
$ws.Range("T30").Value = 0.007438669410940644

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 15.79677433333333
$ws.Range("H31").Value = 47.390323
$ws.Range("I31").Value = 0.01775072507139627
$ws.Range("J31").Value = 0.01833332435500452
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 4.991516
$ws.Range("N31").Value = 9.983032
$ws.Range("O31").Value = 0.1719790502865958
$ws.Range("P31").Value = 0.1465065426562639
$ws.Range("Q31").Value = 78.84985183322266
$ws.Range("R31").Value = 473.099110999336
$ws.Range("S31").Value = 0.003052752839677196
$ws.Range("T31").Value = 0.00268595196664759

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 53.1948375
$ws.Range("H32").Value = 106.389675
$ws.Range("I32").Value = 0.05977466764766092
$ws.Range("J32").Value = 0.0411576941520005
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 13.939895
$ws.Range("N32").Value = 27.87979
$ws.Range("O32").Value = 0.4802889349037177
$ws.Range("P32").Value = 0.4091514124048364
$ws.Range("Q32").Value = 741.5304492920625
$ws.Range("R32").Value = 2966.12179716825
$ws.Range("S32").Value = 0.02870911145871877
$ws.Range("T32").Value = 0.01683972869361728

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 53.1948375
$ws.Range("H33").Value = 106.389675
$ws.Range("I33").Value = 0.05977466764766092
$ws.Range("J33").Value = 0.0411576941520005
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 0.5122363333333334
$ws.Range("N33").Value = 1.536709
$ws.Range("O33").Value = 0.01764872999083942
$ws.Range("P33").Value = 0.02255205859890708
$ws.Range("Q33").Value = 27.2483285132625
$ws.Range("R33").Value = 163.489971079575
$ws.Range("S33").Value = 0.001054946969605732
$ws.Range("T33").Value = 0.0009281907303118108

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 53.1948375
$ws.Range("H34").Value = 106.389675
$ws.Range("I34").Value = 0.05977466764766092
$ws.Range("J34").Value = 0.0411576941520005
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 0.1702606666666667
$ws.Range("N34").Value = 0.5107820000000001
$ws.Range("O34").Value = 0.005866207331499287
$ws.Range("P34").Value = 0.007496009716391951
$ws.Range("Q34").Value = 9.056988495975002
$ws.Range("R34").Value = 54.34193097585
$ws.Range("S34").Value = 0.0003506505935926418
$ws.Range("T34").Value = 0.000308518475267684

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 53.1948375
$ws.Range("H35").Value = 106.389675
$ws.Range("I35").Value = 0.05977466764766092
$ws.Range("J35").Value = 0.0411576941520005
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 0.194159
$ws.Range("N35").Value = 0.5824769999999999
$ws.Range("O35").Value = 0.006689607010093757
$ws.Range("P35").Value = 0.008548173685789306
$ws.Range("Q35").Value = 10.3282564541625
$ws.Range("R35").Value = 61.96953872497499
$ws.Range("S35").Value = 0.000399869035721817
$ws.Range("T35").Value = 0.0003518231181178952

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 53.1948375
$ws.Range("H36").Value = 106.389675
$ws.Range("I36").Value = 0.05977466764766092
$ws.Range("J36").Value = 0.0411576941520005
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 9.215910000000001
$ws.Range("N36").Value = 27.64773
$ws.Range("O36").Value = 0.317527470477254
$ws.Range("P36").Value = 0.4057458029378115
$ws.Range("Q36").Value = 490.2388348646251
$ws.Range("R36").Value = 2941.43300918775
$ws.Range("S36").Value = 0.01898009901678033
$ws.Range("T36").Value = 0.01669956166077231

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 53.1948375
$ws.Range("H37").Value = 106.389675
$ws.Range("I37").Value = 0.05977466764766092
$ws.Range("J37").Value = 0.0411576941520005
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 4.991516
$ws.Range("N37").Value = 9.983032
$ws.Range("O37").Value = 0.1719790502865958
$ws.Range("P37").Value = 0.1465065426562639
$ws.Range("Q37").Value = 265.52288249865
$ws.Range("R37").Value = 1062.0915299946
$ws.Range("S37").Value = 0.01027999057324163
$ws.Range("T37").Value = 0.006029871473913524
